$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values: format cell as Text first so number-looking strings
# (e.g. "1.000", "290.13") are preserved verbatim instead of being
# parsed into numeric values, then clear the format back off so no
# stray style survives the edit.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '22.059.08'
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.554.33'
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '290.13'
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.3934'
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3213'
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '43.69'
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07240'
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '1.072'
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.654'
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '18.70'
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.00001126'
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '1.556.35'
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '6.611'
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.06595'
$cell.ClearFormats()
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '83.31'
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.281'
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '15.45'
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '11.25'
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '22.069.13'
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.372'
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.423'
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '148.54'
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '18.51'
$cell.ClearFormats()
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '4.880'
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.718.86'
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '118.46'
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.9936'
$cell.ClearFormats()
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.789'
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.08284'
$cell.ClearFormats()
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.606'
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '9.004'
$cell.ClearFormats()
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.02257'
$cell.ClearFormats()
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.06054'
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '5.103'
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '1.212'
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.2035'
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '10.64'
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.5798'
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.744'
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '12.87'
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.5573'
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.891'
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '117.51'
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.130'
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.06818'
$cell.ClearFormats()

# Column B, C, E values: plain text, safe to assign directly.
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("E7").Value = '  +3.61%  '
$ws.Range("E8").Value = '  -2.04%  '
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("E10").Value = '  -1.74%  '
$ws.Range("E11").Value = '  -5.59%  '
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("E14").Value = '  -7.52%  '
$ws.Range("E15").Value = '  +4.73%  '
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("E19").Value = '  -3.66%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("E22").Value = '  -4.57%  '
$ws.Range("E23").Value = '  -3.96%  '
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("E25").Value = '  +3.75%  '
$ws.Range("E26").Value = '  -5.57%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  -4.17%  '
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("E31").Value = '  -3.63%  '
$ws.Range("E32").Value = '  -8.27%  '
$ws.Range("E33").Value = '  -2.29%  '
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("E35").Value = '  -16.52%  '
$ws.Range("E36").Value = '  -4.26%  '
$ws.Range("E37").Value = '  -4.04%  '
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("E39").Value = '  -4.81%  '
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("E41").Value = '  -5.71%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  -3.60%  '
$ws.Range("E44").Value = '  -4.53%  '
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("E46").Value = '  -6.80%  '
$ws.Range("E47").Value = '  -5.50%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E48").Value = '  -4.21%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E49").Value = '  -4.61%  '
$ws.Range("E50").Value = '  -4.08%  '
$ws.Range("E51").Value = '  -3.58%  '
